# Update script matching commit "Atualizado por script em 19-12-2023 18:56"
# - Reorders a few existing match rows (the scraper re-sorted same-kickoff
#   matches), and appends 19 new match rows (80-98) scraped afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Part 1: rows whose B:V content changes in place (A / Indice stays).
# ------------------------------------------------------------------
$updates = @(
    @{Row=50; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Zob Ahan"; G=1; H="Tractor"; I=0;
      J=3.1; K="26/10/2023 07:42"; L=3.13; M="27/10/2023 16:17";
      N=2.74; O="26/10/2023 07:42"; P=2.69; Q="27/10/2023 16:18";
      R=2.35; S="26/10/2023 07:42"; T=2.63; U="27/10/2023 16:18";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/zob-ahan-tractor/I9W8sTvN/"},

    @{Row=51; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Esteghlal F.C."; G=1; H="Aluminium Arak"; I=0;
      J=1.49; K="26/10/2023 07:42"; L=1.54; M="27/10/2023 16:24";
      N=3.49; O="26/10/2023 07:42"; P=3.42; Q="27/10/2023 16:25";
      R=6.44; S="26/10/2023 07:42"; T=7.84; U="27/10/2023 16:25";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/esteghlal-teh-aluminium-arak/lzS4r9PG/"},

    @{Row=66; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Malavan"; G=1; H="Esteghlal Khuzestan"; I=1;
      J=1.81; K="09/11/2023 00:42"; L=1.39; M="10/11/2023 12:29";
      N=2.92; O="09/11/2023 00:42"; P=3.45; Q="10/11/2023 12:29";
      R=4.57; S="09/11/2023 00:42"; T=7.84; U="10/11/2023 12:29";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/malavan-esteghlal-khuzestan/CfSExSU9/"},

    @{Row=67; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Gol Gohar"; G=2; H="Mes Rafsanjan"; I=0;
      J=2.54; K="09/11/2023 00:42"; L=2.89; M="10/11/2023 12:26";
      N=2.66; O="09/11/2023 00:42"; P=2.45; Q="10/11/2023 12:26";
      R=2.9; S="09/11/2023 00:42"; T=3.16; U="10/11/2023 12:26";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/gol-gohar-mes-rafsanjan/fJHJynpG/"},

    @{Row=74; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Mes Rafsanjan"; G=0; H="Sanat Naft"; I=2;
      J=1.53; K="22/11/2023 15:12"; L=1.58; M="24/11/2023 12:25";
      N=3.43; O="22/11/2023 15:12"; P=3.36; Q="24/11/2023 12:25";
      R=6.31; S="22/11/2023 15:12"; T=7.1; U="24/11/2023 12:25";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-sanat-naft/MVnvGh0r/"},

    @{Row=75; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Aluminium Arak"; G=1; H="Paykan"; I=1;
      J=1.85; K="22/11/2023 15:12"; L=2.1; M="24/11/2023 12:10";
      N=2.8; O="22/11/2023 15:12"; P=2.36; Q="24/11/2023 12:10";
      R=4.58; S="22/11/2023 15:12"; T=4.43; U="24/11/2023 12:10";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-paykan/2o67LExL/"},

    @{Row=76; B="iran"; C="persian-gulf-pro-league"; D="2023-2024";
      F="Havadar SC"; G=0; H="Malavan"; I=1;
      J=2.51; K="22/11/2023 15:12"; L=3.16; M="24/11/2023 11:20";
      N=2.65; O="22/11/2023 15:12"; P=2.5; Q="24/11/2023 11:20";
      R=2.96; S="22/11/2023 15:12"; T=2.81; U="24/11/2023 11:33";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-malavan/tx5BKYiR/"}
)

foreach ($r in $updates) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H
    $ws.Range("I" + $r.Row).Value = $r.I
    $ws.Range("J" + $r.Row).Value = $r.J
    $ws.Range("K" + $r.Row).Value = $r.K
    $ws.Range("L" + $r.Row).Value = $r.L
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("N" + $r.Row).Value = $r.N
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("P" + $r.Row).Value = $r.P
    $ws.Range("Q" + $r.Row).Value = $r.Q
    $ws.Range("R" + $r.Row).Value = $r.R
    $ws.Range("S" + $r.Row).Value = $r.S
    $ws.Range("T" + $r.Row).Value = $r.T
    $ws.Range("U" + $r.Row).Value = $r.U
    $ws.Range("V" + $r.Row).Value = $r.V
}

# ------------------------------------------------------------------
# Part 2: append 19 brand-new rows (80-98), copying the formatting
# of the final pre-existing row (79) so styles (borders/number
# formats) carry over, then filling in the values.
# ------------------------------------------------------------------
$lastDataRow = 79
$newRowCount = 19
$firstNewRow = $lastDataRow + 1
$lastNewRow = $lastDataRow + $newRowCount

$ws.Range("A" + $lastDataRow + ":V" + $lastDataRow).Copy()
$ws.Range("A" + $firstNewRow + ":V" + $lastNewRow).PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{Row=80; A=79; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45267.52083333334;
      F="Gol Gohar"; G=0; H="Tractor"; I=1;
      J=3.11; K="06/12/2023 00:42"; L=3.41; M="07/12/2023 12:08";
      N=2.77; O="06/12/2023 00:42"; P=2.8; Q="07/12/2023 12:08";
      R=2.36; S="06/12/2023 00:42"; T=2.37; U="07/12/2023 12:08";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/gol-gohar-tractor/U1waBUED/"},

    @{Row=81; A=80; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45268.52083333334;
      F="Malavan"; G=0; H="Shams Azar Qazvin"; I=2;
      J=2.39; K="07/12/2023 00:42"; L=3; M="08/12/2023 12:09";
      N=2.66; O="07/12/2023 00:42"; P=2.42; Q="08/12/2023 11:55";
      R=3.2; S="07/12/2023 00:42"; T=3.08; U="08/12/2023 12:09";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/malavan-shams-azar-qazvin/M9upZX6K/"},

    @{Row=82; A=81; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45268.52083333334;
      F="Paykan"; G=0; H="Zob Ahan"; I=0;
      J=3.54; K="07/12/2023 00:42"; L=4.33; M="08/12/2023 12:13";
      N=2.57; O="07/12/2023 00:42"; P=2.37; Q="08/12/2023 12:13";
      R=2.3; S="07/12/2023 00:42"; T=2.39; U="08/12/2023 12:13";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/paykan-zob-ahan/4YUlYiMQ/"},

    @{Row=83; A=82; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45268.52083333334;
      F="Sanat Naft"; G=1; H="Aluminium Arak"; I=2;
      J=2.84; K="07/12/2023 00:42"; L=2.72; M="08/12/2023 12:29";
      N=2.59; O="07/12/2023 00:42"; P=2.56; Q="08/12/2023 12:29";
      R=2.73; S="07/12/2023 00:42"; T=3.19; U="08/12/2023 12:29";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sanat-naft-aluminium-arak/rJttzEiE/"},

    @{Row=84; A=83; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45269.52083333334;
      F="Nassaji Mazandaran"; G=1; H="Esteghlal Khuzestan"; I=0;
      J=1.84; K="08/12/2023 00:43"; L=2; M="09/12/2023 12:25";
      N=2.89; O="08/12/2023 00:43"; P=2.68; Q="09/12/2023 12:25";
      R=4.6; S="08/12/2023 00:43"; T=5.06; U="09/12/2023 11:53";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mazandaran-esteghlal-khuzestan/tfMGTBUs/"},

    @{Row=85; A=84; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45269.52083333334;
      F="Sepahan"; G=4; H="Mes Rafsanjan"; I=1;
      J=1.54; K="08/12/2023 00:43"; L=1.6; M="09/12/2023 12:17";
      N=3.43; O="08/12/2023 00:43"; P=3.46; Q="09/12/2023 12:17";
      R=5.99; S="08/12/2023 00:43"; T=6.48; U="09/12/2023 12:17";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sepahan-mes-rafsanjan/zHLKSVql/"},

    @{Row=86; A=85; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45269.54166666666;
      F="Persepolis"; G=1; H="Havadar SC"; I=0;
      J=1.31; K="08/12/2023 01:12"; L=1.28; M="09/12/2023 10:17";
      N=4.34; O="08/12/2023 01:12"; P=4.85; Q="09/12/2023 12:29";
      R=9.31; S="08/12/2023 01:12"; T=12.37; U="09/12/2023 12:29";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/persepolis-havadar-sc/E7KORkaf/"},

    @{Row=87; A=86; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45269.63541666666;
      F="Foolad"; G=0; H="Esteghlal F.C."; I=0;
      J=4.28; K="08/12/2023 03:42"; L=5.61; M="09/12/2023 15:06";
      N=2.82; O="08/12/2023 03:42"; P=3.05; Q="09/12/2023 15:06";
      R=1.93; S="08/12/2023 03:42"; T=1.77; U="09/12/2023 15:06";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/foolad-esteghlal-teh/WrFTQ9E0/"},

    @{Row=88; A=87; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45273.52083333334;
      F="Sanat Naft"; G=0; H="Gol Gohar"; I=1;
      J=3.03; K="12/12/2023 00:42"; L=3.95; M="13/12/2023 12:26";
      N=2.75; O="12/12/2023 00:42"; P=2.67; Q="13/12/2023 12:26";
      R=2.42; S="12/12/2023 00:42"; T=2.26; U="13/12/2023 12:26";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sanat-naft-gol-gohar/6FhzIKMm/"},

    @{Row=89; A=88; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45274.52083333334;
      F="Sepahan"; G=3; H="Foolad"; I=1;
      J=1.4; K="30/09/2023 13:43"; L=1.38; M="14/12/2023 12:28";
      N=3.91; O="30/09/2023 13:43"; P=4.19; Q="14/12/2023 12:28";
      R=6.8; S="30/09/2023 13:43"; T=9.25; U="14/12/2023 12:28";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sepahan-foolad/lvbVI07s/"},

    @{Row=90; A=89; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45274.57291666666;
      F="Persepolis"; G=1; H="Esteghlal F.C."; I=1;
      J=1.96; K="30/09/2023 13:43"; L=2.4; M="14/12/2023 13:40";
      N=2.88; O="30/09/2023 13:43"; P=2.62; Q="14/12/2023 13:40";
      R=3.87; S="30/09/2023 13:43"; T=3.68; U="14/12/2023 13:40";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/persepolis-esteghlal-teh/vTjrGbi0/"},

    @{Row=91; A=90; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45274.625;
      F="Shams Azar Qazvin"; G=2; H="Nassaji Mazandaran"; I=1;
      J=1.91; K="13/12/2023 03:12"; L=2.33; M="14/12/2023 14:59";
      N=2.96; O="13/12/2023 03:12"; P=2.78; Q="14/12/2023 13:36";
      R=4.11; S="13/12/2023 03:12"; T=3.55; U="14/12/2023 14:59";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-mazandaran/x0UC0JiD/"},

    @{Row=92; A=91; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45278.52083333334;
      F="Mes Rafsanjan"; G=1; H="Malavan"; I=4;
      J=2; K="17/12/2023 00:42"; L=2.38; M="18/12/2023 12:29";
      N=2.81; O="17/12/2023 00:42"; P=2.64; Q="18/12/2023 12:29";
      R=3.98; S="17/12/2023 00:42"; T=3.06; U="18/12/2023 12:29";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-malavan/4dDyPmqD/"},

    @{Row=93; A=92; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45278.52083333334;
      F="Tractor"; G=2; H="Paykan"; I=0;
      J=1.47; K="17/12/2023 00:42"; L=1.52; M="18/12/2023 11:50";
      N=3.59; O="17/12/2023 00:42"; P=3.64; Q="18/12/2023 11:50";
      R=6.76; S="17/12/2023 00:42"; T=7.32; U="18/12/2023 11:50";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/tractor-paykan/AZDXPTT6/"},

    @{Row=94; A=93; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45279.47916666666;
      F="Havadar SC"; G=0; H="Sepahan"; I=5;
      J=6.88; K="17/12/2023 23:42"; L=7.62; M="19/12/2023 11:28";
      N=3.78; O="17/12/2023 23:42"; P=3.69; Q="19/12/2023 11:28";
      R=1.44; S="17/12/2023 23:42"; T=1.5; U="19/12/2023 11:26";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-sepahan/GWGuO7bJ/"},

    @{Row=95; A=94; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45279.52083333334;
      F="Aluminium Arak"; G=0; H="Gol Gohar"; I=0;
      J=2.74; K="18/12/2023 00:42"; L=3.26; M="19/12/2023 12:23";
      N=2.56; O="18/12/2023 00:42"; P=2.29; Q="19/12/2023 12:23";
      R=2.86; S="18/12/2023 00:42"; T=3.04; U="19/12/2023 12:23";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-gol-gohar/betQVJHP/"},

    @{Row=96; A=95; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45279.52083333334;
      F="Esteghlal Khuzestan"; G=2; H="Persepolis"; I=2;
      J=12.94; K="18/12/2023 01:12"; L=10.21; M="19/12/2023 12:29";
      N=4.63; O="18/12/2023 01:12"; P=4; Q="19/12/2023 12:29";
      R=1.24; S="18/12/2023 01:12"; T=1.39; U="19/12/2023 12:29";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/esteghlal-khuzestan-persepolis/jmYiRwQt/"},

    @{Row=97; A=96; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45279.52083333334;
      F="Shams Azar Qazvin"; G=0; H="Sanat Naft"; I=0;
      J=2.01; K="18/12/2023 00:42"; L=1.85; M="19/12/2023 10:32";
      N=2.9; O="18/12/2023 00:42"; P=2.84; Q="19/12/2023 10:32";
      R=3.81; S="18/12/2023 00:42"; T=5.56; U="19/12/2023 11:41";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-sanat-naft/K4sMWa2J/"},

    @{Row=98; A=97; B="iran"; C="persian-gulf-pro-league"; D="2023-2024"; E=45279.58333333334;
      F="Zob Ahan"; G=0; H="Foolad"; I=0;
      J=2.23; K="18/12/2023 02:12"; L=2.07; M="19/12/2023 13:58";
      N=2.6; O="18/12/2023 02:12"; P=2.58; Q="19/12/2023 13:58";
      R=3.66; S="18/12/2023 02:12"; T=4.98; U="19/12/2023 13:58";
      V="https://www.betexplorer.com/football/iran/persian-gulf-pro-league/zob-ahan-foolad/UwXeQcum/"}
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H
    $ws.Range("I" + $r.Row).Value = $r.I
    $ws.Range("J" + $r.Row).Value = $r.J
    $ws.Range("K" + $r.Row).Value = $r.K
    $ws.Range("L" + $r.Row).Value = $r.L
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("N" + $r.Row).Value = $r.N
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("P" + $r.Row).Value = $r.P
    $ws.Range("Q" + $r.Row).Value = $r.Q
    $ws.Range("R" + $r.Row).Value = $r.R
    $ws.Range("S" + $r.Row).Value = $r.S
    $ws.Range("T" + $r.Row).Value = $r.T
    $ws.Range("U" + $r.Row).Value = $r.U
    $ws.Range("V" + $r.Row).Value = $r.V
}

Write-Host "Edit complete"
